$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuites")

$ws.Range("A10").Value = "OrderTestCases"
$ws.Range("C10").Value = "Y"

$ws.Range("E11").Select()
